$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value2 = 22736682
$ws.Cells.Item(62, 9).Value2 = 9058.117
$ws.Cells.Item(62, 11).Value2 = 9058.117
$ws.Cells.Item(62, 13).Value2 = -8434.117

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value2 = 22736682
$ws.Cells.Item(65, 9).Value2 = 9058.117
$ws.Cells.Item(65, 11).Value2 = 45290.585
$ws.Cells.Item(65, 13).Value2 = -42170.585

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value2 = 15785.421
$ws.Cells.Item(70, 9).Value2 = 27336.875
$ws.Cells.Item(70, 10).Value2 = 12705.033
$ws.Cells.Item(70, 11).Value2 = 82010.625
$ws.Cells.Item(70, 12).Value2 = 38115.099
$ws.Cells.Item(70, 13).Value2 = -81740.625
$ws.Cells.Item(70, 14).Value2 = -38655.099

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value2 = 15785.421
$ws.Cells.Item(73, 9).Value2 = 27336.875
$ws.Cells.Item(73, 10).Value2 = 12705.033
$ws.Cells.Item(73, 11).Value2 = 82010.625
$ws.Cells.Item(73, 12).Value2 = 38115.099
$ws.Cells.Item(73, 13).Value2 = -81074.625
$ws.Cells.Item(73, 14).Value2 = -39987.099

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value2 = 1396.7391
$ws.Cells.Item(92, 9).Value2 = 530.5238000000001
$ws.Cells.Item(92, 11).Value2 = 530.5238000000001
$ws.Cells.Item(92, 13).Value2 = 717.4761999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value2 = 3928898.2
$ws.Cells.Item(112, 10).Value2 = 5439694.5
$ws.Cells.Item(112, 12).Value2 = 16319083.5
$ws.Cells.Item(112, 14).Value2 = -16321299.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value2 = 2876.23
$ws.Cells.Item(138, 9).Value2 = 1883.4
$ws.Cells.Item(138, 10).Value2 = 3538.1167
$ws.Cells.Item(138, 11).Value2 = 5650.200000000001
$ws.Cells.Item(138, 12).Value2 = 10614.3501
$ws.Cells.Item(138, 13).Value2 = -510.2000000000007
$ws.Cells.Item(138, 14).Value2 = -20894.3501

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value2 = 1494.28
$ws.Cells.Item(45, 9).Value2 = 982.3125
$ws.Cells.Item(45, 11).Value2 = 982.3125
$ws.Cells.Item(45, 13).Value2 = -605.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value2 = 4883.4736
$ws.Cells.Item(61, 9).Value2 = 2275.3333
$ws.Cells.Item(61, 10).Value2 = 7230.8
$ws.Cells.Item(61, 11).Value2 = 2275.3333
$ws.Cells.Item(61, 12).Value2 = 7230.8
$ws.Cells.Item(61, 13).Value2 = -2063.3333
$ws.Cells.Item(61, 14).Value2 = -7654.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value2 = 3030.6667
$ws.Cells.Item(63, 9).Value2 = 2974.2354
$ws.Cells.Item(63, 11).Value2 = 2974.2354
$ws.Cells.Item(63, 13).Value2 = -2288.2354

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value2 = 3030.6667
$ws.Cells.Item(66, 9).Value2 = 2974.2354
$ws.Cells.Item(66, 11).Value2 = 14871.177
$ws.Cells.Item(66, 13).Value2 = -11439.177

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(119, 8).Value2 = 0
$ws.Cells.Item(119, 10).Value2 = 0
$ws.Cells.Item(119, 12).Value2 = 0
$ws.Cells.Item(119, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value2 = 5047.35
$ws.Cells.Item(122, 9).Value2 = 3712.0264
$ws.Cells.Item(122, 10).Value2 = 30418.5
$ws.Cells.Item(122, 11).Value2 = 11136.0792
$ws.Cells.Item(122, 12).Value2 = 91255.5
$ws.Cells.Item(122, 13).Value2 = -8686.0792
$ws.Cells.Item(122, 14).Value2 = -96155.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(124, 8).Value2 = 0
$ws.Cells.Item(124, 10).Value2 = 0
$ws.Cells.Item(124, 12).Value2 = 0
$ws.Cells.Item(124, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(129, 8).Value2 = 0
$ws.Cells.Item(129, 9).Value2 = 0
$ws.Cells.Item(129, 10).Value2 = 0
$ws.Cells.Item(129, 11).Value2 = 0
$ws.Cells.Item(129, 12).Value2 = 0
$ws.Cells.Item(129, 13).ClearContents()
$ws.Cells.Item(129, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value2 = 4883.4736
$ws.Cells.Item(136, 9).Value2 = 2275.3333
$ws.Cells.Item(136, 10).Value2 = 7230.8
$ws.Cells.Item(136, 11).Value2 = 6825.999899999999
$ws.Cells.Item(136, 12).Value2 = 21692.4
$ws.Cells.Item(136, 13).Value2 = -4275.999899999999
$ws.Cells.Item(136, 14).Value2 = -26792.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value2 = 938.05554
$ws.Cells.Item(80, 9).Value2 = 1157.7273
$ws.Cells.Item(80, 10).Value2 = 592.8570999999999
$ws.Cells.Item(80, 11).Value2 = 1157.7273
$ws.Cells.Item(80, 12).Value2 = 592.8570999999999
$ws.Cells.Item(80, 13).Value2 = -159.7273
$ws.Cells.Item(80, 14).Value2 = -2588.8571

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value2 = 938.05554
$ws.Cells.Item(83, 9).Value2 = 1157.7273
$ws.Cells.Item(83, 10).Value2 = 592.8570999999999
$ws.Cells.Item(83, 11).Value2 = 5788.636500000001
$ws.Cells.Item(83, 12).Value2 = 2964.2855
$ws.Cells.Item(83, 13).Value2 = -796.6365000000005
$ws.Cells.Item(83, 14).Value2 = -12948.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 20851902
$ws.Cells.Item(31, 9).Value2 = 1493.2258
$ws.Cells.Item(31, 11).Value2 = 1493.2258
$ws.Cells.Item(31, 13).Value2 = -1198.2258

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value2 = 20851902
$ws.Cells.Item(34, 9).Value2 = 1493.2258
$ws.Cells.Item(34, 11).Value2 = 1493.2258
$ws.Cells.Item(34, 13).Value2 = -1291.2258

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(39, 8).Value2 = 15867.25
$ws.Cells.Item(39, 9).Value2 = 4489.6665
$ws.Cells.Item(39, 10).Value2 = 50000
$ws.Cells.Item(39, 11).Value2 = 4489.6665
$ws.Cells.Item(39, 12).Value2 = 50000
$ws.Cells.Item(39, 13).Value2 = -4098.6665
$ws.Cells.Item(39, 14).Value2 = -50782

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(49, 8).Value2 = 15867.25
$ws.Cells.Item(49, 9).Value2 = 4489.6665
$ws.Cells.Item(49, 10).Value2 = 50000
$ws.Cells.Item(49, 11).Value2 = 4489.6665
$ws.Cells.Item(49, 12).Value2 = 50000
$ws.Cells.Item(49, 13).Value2 = -4307.6665
$ws.Cells.Item(49, 14).Value2 = -50364

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(141, 8).Value2 = 76218.836
$ws.Cells.Item(141, 10).Value2 = 77304.46000000001
$ws.Cells.Item(141, 12).Value2 = 77304.46000000001
$ws.Cells.Item(141, 14).Value2 = -87664.46000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value2 = 1457.1875
$ws.Cells.Item(102, 9).Value2 = 1535.7046
$ws.Cells.Item(102, 11).Value2 = 1535.7046
$ws.Cells.Item(102, 13).Value2 = 86.29539999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value2 = 1417.1364
$ws.Cells.Item(122, 9).Value2 = 1208.9
$ws.Cells.Item(122, 10).Value2 = 3499.5
$ws.Cells.Item(122, 11).Value2 = 3626.7
$ws.Cells.Item(122, 12).Value2 = 10498.5
$ws.Cells.Item(122, 13).Value2 = -1176.7
$ws.Cells.Item(122, 14).Value2 = -15398.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value2 = 8946.037
$ws.Cells.Item(126, 10).Value2 = 8537.933999999999
$ws.Cells.Item(126, 12).Value2 = 25613.802
$ws.Cells.Item(126, 14).Value2 = -30553.802

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value2 = 22610.8
$ws.Cells.Item(132, 9).Value2 = 23259.717
$ws.Cells.Item(132, 11).Value2 = 69779.151
$ws.Cells.Item(132, 13).Value2 = -67249.151

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(134, 8).Value2 = 27166
$ws.Cells.Item(134, 10).Value2 = 27166
$ws.Cells.Item(134, 12).Value2 = 81498
$ws.Cells.Item(134, 14).Value2 = -86568

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value2 = 30304.834
$ws.Cells.Item(136, 10).Value2 = 30304.834
$ws.Cells.Item(136, 12).Value2 = 90914.50199999999
$ws.Cells.Item(136, 14).Value2 = -96014.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value2 = 1743.55
$ws.Cells.Item(68, 9).Value2 = 1722.9429
$ws.Cells.Item(68, 10).Value2 = 1887.8
$ws.Cells.Item(68, 11).Value2 = 1722.9429
$ws.Cells.Item(68, 12).Value2 = 1887.8
$ws.Cells.Item(68, 13).Value2 = -973.9429
$ws.Cells.Item(68, 14).Value2 = -3385.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value2 = 1743.55
$ws.Cells.Item(71, 9).Value2 = 1722.9429
$ws.Cells.Item(71, 10).Value2 = 1887.8
$ws.Cells.Item(71, 11).Value2 = 8614.7145
$ws.Cells.Item(71, 12).Value2 = 9439
$ws.Cells.Item(71, 13).Value2 = -4870.7145
$ws.Cells.Item(71, 14).Value2 = -16927

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value2 = 1898.5
$ws.Cells.Item(82, 10).Value2 = 1926.5
$ws.Cells.Item(82, 12).Value2 = 1926.5
$ws.Cells.Item(82, 14).Value2 = -2648.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value2 = 1898.5
$ws.Cells.Item(85, 10).Value2 = 1926.5
$ws.Cells.Item(85, 12).Value2 = 1926.5
$ws.Cells.Item(85, 14).Value2 = -4422.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value2 = 2260.566
$ws.Cells.Item(136, 9).Value2 = 1855.3954
$ws.Cells.Item(136, 11).Value2 = 5566.1862
$ws.Cells.Item(136, 13).Value2 = -3016.1862

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value2 = 10303
$ws.Cells.Item(62, 9).Value2 = 8694.5
$ws.Cells.Item(62, 11).Value2 = 8694.5
$ws.Cells.Item(62, 13).Value2 = -8070.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value2 = 10303
$ws.Cells.Item(65, 9).Value2 = 8694.5
$ws.Cells.Item(65, 11).Value2 = 43472.5
$ws.Cells.Item(65, 13).Value2 = -40352.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value2 = 54333
$ws.Cells.Item(135, 10).Value2 = 56124.625
$ws.Cells.Item(135, 12).Value2 = 56124.625
$ws.Cells.Item(135, 14).Value2 = -66264.625
